$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2,2).Value = 24
$ws.Cells.Item(2,3).Value = "Walker Kessler"
$ws.Cells.Item(2,4).Value = "C"
$ws.Cells.Item(2,5).Value = "7-1"
$ws.Cells.Item(2,6).Value = 245
$ws.Cells.Item(2,7).Value = "July 26, 2001"
$ws.Cells.Item(2,9).Value = "R"
$ws.Cells.Item(2,10).Value = "UNC, Auburn"
$ws.Cells.Item(2,11).Value = "https://www.basketball-reference.com/players/k/kesslwa01.html"
$ws.Cells.Item(3,2).Value = 0
$ws.Cells.Item(3,3).Value = "Jordan Clarkson"
$ws.Cells.Item(3,4).Value = "SG"
$ws.Cells.Item(3,5).Value = "6-4"
$ws.Cells.Item(3,6).Value = 194
$ws.Cells.Item(3,7).Value = "June 7, 1992"
$ws.Cells.Item(3,9).Value = "8"
$ws.Cells.Item(3,10).Value = "Tulsa, Missouri"
$ws.Cells.Item(3,11).Value = "https://www.basketball-reference.com/players/c/clarkjo01.html"
$ws.Cells.Item(6,2).Value = 0
$ws.Cells.Item(6,3).Value = "Talen Horton-Tucker"
$ws.Cells.Item(6,5).Value = "6-4"
$ws.Cells.Item(6,6).Value = 234
$ws.Cells.Item(6,7).Value = "November 25, 2000"
$ws.Cells.Item(6,9).Value = "3"
$ws.Cells.Item(6,10).Value = "Iowa State"
$ws.Cells.Item(6,11).Value = "https://www.basketball-reference.com/players/h/hortota01.html"
$ws.Cells.Item(7,2).Value = 22
$ws.Cells.Item(7,3).Value = "Rudy Gay"
$ws.Cells.Item(7,4).Value = "PF"
$ws.Cells.Item(7,5).Value = "6-8"
$ws.Cells.Item(7,6).Value = 250
$ws.Cells.Item(7,7).Value = "August 17, 1986"
$ws.Cells.Item(7,9).Value = "16"
$ws.Cells.Item(7,10).Value = "UConn"
$ws.Cells.Item(7,11).Value = "https://www.basketball-reference.com/players/g/gayru01.html"
$ws.Cells.Item(8,2).Value = 2
$ws.Cells.Item(8,3).Value = "Collin Sexton"
$ws.Cells.Item(8,4).Value = "SG"
$ws.Cells.Item(8,5).Value = "6-1"
$ws.Cells.Item(8,6).Value = 190
$ws.Cells.Item(8,7).Value = "January 4, 1999"
$ws.Cells.Item(8,9).Value = "4"
$ws.Cells.Item(8,10).Value = "Alabama"
$ws.Cells.Item(8,11).Value = "https://www.basketball-reference.com/players/s/sextoco01.html"
$ws.Cells.Item(12,2).Value = 11
$ws.Cells.Item(12,3).Value = "Kris Dunn"
$ws.Cells.Item(12,4).Value = "PG"
$ws.Cells.Item(12,5).Value = "6-3"
$ws.Cells.Item(12,6).Value = 205
$ws.Cells.Item(12,7).Value = "March 18, 1994"
$ws.Cells.Item(12,9).Value = "6"
$ws.Cells.Item(12,10).Value = "Providence"
$ws.Cells.Item(12,11).Value = "https://www.basketball-reference.com/players/d/dunnkr01.html"
$ws.Cells.Item(15,2).Value = 33
$ws.Cells.Item(15,3).Value = "Johnny Juzang (TW)"
$ws.Cells.Item(15,4).Value = "SF"
$ws.Cells.Item(15,5).Value = "6-7"
$ws.Cells.Item(15,6).Value = 215
$ws.Cells.Item(15,7).Value = "March 17, 2001"
$ws.Cells.Item(15,9).Value = "R"
$ws.Cells.Item(15,10).Value = "Kentucky, UCLA"
$ws.Cells.Item(15,11).Value = "https://www.basketball-reference.com/players/j/juzanjo01.html"
$ws.Cells.Item(16,2).Value = 25
$ws.Cells.Item(16,3).Value = "Micah Potter (TW)"
$ws.Cells.Item(16,4).Value = "PF"
$ws.Cells.Item(16,5).Value = "6-10"
$ws.Cells.Item(16,6).Value = 248
$ws.Cells.Item(16,7).Value = "April 6, 1998"
$ws.Cells.Item(16,9).Value = "1"
$ws.Cells.Item(16,10).Value = "Ohio State, Wisconsin"
$ws.Cells.Item(16,11).Value = "https://www.basketball-reference.com/players/p/pottemi01.html"
